$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing column (K) into the new column (L)
# for the border-only row (3), the year-header row (4) and the data row (5).
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("L4").Value = 2021

$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("L5").Value = 269

$excel.CutCopyMode = $false

# Update the active selection to match the target workbook state.
$ws.Range("N3").Select()
